$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.870.44"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").Value = "3.933.11"
$ws.Range("E3").Value = "  +4.62%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'605.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.52%  "
$ws.Range("D6").Value = "'166.22"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("D7").Value = "3.931.75"
$ws.Range("E7").Value = "  +4.67%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'0.532"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.25%  "
$ws.Range("E10").Value = "  -3.29%  "
$ws.Range("D11").Value = "'6.43"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.29%  "
$ws.Range("D12").Value = "'0.463"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.92%  "
$ws.Range("D13").Value = "'37.37"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.08%  "
$ws.Range("D14").Value = "'0.0000247"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.60%  "
$ws.Range("D15").Value = "4.582.58"
$ws.Range("E15").Value = "  +4.38%  "
$ws.Range("D16").Value = "3.924.83"
$ws.Range("E16").Value = "  +4.48%  "
$ws.Range("D17").Value = "68.965.30"
$ws.Range("E17").Value = "  -0.18%  "
$ws.Range("D18").Value = "'7.50"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.84%  "
$ws.Range("E19").Value = "  -0.60%  "
$ws.Range("D20").Value = "'17.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.13%  "
$ws.Range("D21").Value = "'11.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.68%  "
$ws.Range("D22").Value = "'489.04"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.34%  "
$ws.Range("D23").Value = "'0.725"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").Value = "'0.0000166"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +11.22%  "
$ws.Range("D25").Value = "'84.53"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.39%  "
$ws.Range("D26").Value = "'2.28"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.26%  "
$ws.Range("D27").Value = "'12.15"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.24%  "
$ws.Range("D28").Value = "'10.19"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.99%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("E30").Value = "  -0.40%  "
$ws.Range("D31").Value = "4.077.83"
$ws.Range("E31").Value = "  +4.37%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "'2.40"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.50%  "
$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").Value = "'7.90"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.62%  "
$ws.Range("D34").Value = "'32.38"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.95%  "
$ws.Range("D35").Value = "3.876.74"
$ws.Range("E35").Value = "  +4.84%  "
$ws.Range("D36").Value = "'0.108"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("E37").Value = "  +3.27%  "
$ws.Range("E38").Value = "  +1.64%  "
$ws.Range("D39").Value = "'5.95"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.12%  "
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("D41").Value = "'0.323"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.52%  "
$ws.Range("D42").Value = "'442.07"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.78%  "
$ws.Range("D43").Value = "'3.01"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.65%  "
$ws.Range("D44").Value = "'2.01"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("E45").Value = "  -0.22%  "
$ws.Range("D46").Value = "'8.53"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.95%  "
$ws.Range("D48").Value = "2.852.60"
$ws.Range("E48").Value = "  +2.07%  "
$ws.Range("D49").Value = "'26.39"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +11.96%  "
$ws.Range("D50").Value = "'141.76"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.16%  "
$ws.Range("D51").Value = "'0.0358"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.31%  "
